$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estudos")

# Row 45 (date 44870): fill in HORA F / DESCANSO that had been left blank,
# which also recalculates DIF (D45) and ÚTEIS (F45) via their table formulas,
# and set ASSUNTO (G45) to "HARD".
$ws.Range("C45").Value = 0.66666666666666663   # 16:00 -> HORA F
$ws.Range("E45").Value = 0.0069444444444444441 # 0:10  -> DESCANSO
$ws.Range("G45").Value = "HARD"                # ASSUNTO

# Row 46 (date 44871): fill in HORA I, ASSUNTO and PRODUÇÃO.
$ws.Range("B46").Value = 0.40277777777777773   # 9:40 -> HORA I
$ws.Range("G46").Value = "HARD"                # ASSUNTO
$ws.Range("H46").Value = "Resolução de atividades banco de dados" # PRODUÇÃO

# The same PRODUÇÃO text also replaces row 45's existing entry (shared string).
$ws.Range("H45").Value = "Resolução de atividades banco de dados"

# Move/update the active selection to reflect where the user ended up.
$ws.Range("H52").Select()
